$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.654.29"
$ws.Cells.Item(2, 5).Value = "  +0.14%  "

$ws.Cells.Item(3, 4).Value = "2.648.81"
$ws.Cells.Item(3, 5).Value = "  -0.05%  "

$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "603.99"
$ws.Cells.Item(5, 5).Value = "  +2.20%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "146.74"
$ws.Cells.Item(6, 5).Value = "  +1.82%  "

$ws.Cells.Item(7, 5).Value = "  +0.03%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.589"
$ws.Cells.Item(8, 5).Value = "  +0.49%  "

$ws.Cells.Item(9, 4).Value = "0.109"
$ws.Cells.Item(9, 5).Value = "  +1.09%  "

$ws.Cells.Item(10, 5).Value = "  -0.35%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.370"

$ws.Cells.Item(12, 5).Value = "  -0.20%  "

$ws.Cells.Item(13, 4).Value = "27.59"
$ws.Cells.Item(13, 5).Value = "  +0.49%  "

$ws.Cells.Item(14, 4).Value = "3.128.09"
$ws.Cells.Item(14, 5).Value = "  +0.15%  "

$ws.Cells.Item(15, 4).Value = "63.502.73"
$ws.Cells.Item(15, 5).Value = "  +0.05%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.0000147"
$ws.Cells.Item(16, 5).Value = "  +0.89%  "

$ws.Cells.Item(17, 4).Value = "2.642.55"
$ws.Cells.Item(17, 5).Value = "  +0.14%  "

$ws.Cells.Item(18, 4).Value = "11.52"
$ws.Cells.Item(18, 5).Value = "  +0.71%  "

$ws.Cells.Item(19, 5).Value = "  +4.37%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "343.53"
$ws.Cells.Item(20, 5).Value = "  +0.85%  "

$ws.Cells.Item(21, 5).Value = "  +2.96%  "

$ws.Cells.Item(22, 5).Value = "  -0.03%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.58"
$ws.Cells.Item(23, 5).Value = "  -3.29%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "66.87"
$ws.Cells.Item(24, 5).Value = "  -0.48%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.70"
$ws.Cells.Item(25, 5).Value = "  +1.21%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "9.06"
$ws.Cells.Item(26, 5).Value = "  +7.61%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "573.81"
$ws.Cells.Item(27, 5).Value = "  +5.00%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.56"
$ws.Cells.Item(28, 5).Value = "  +0.74%  "

$ws.Cells.Item(29, 5).Value = "  -1.59%  "

$ws.Cells.Item(30, 5).Value = "  +2.81%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.00"
$ws.Cells.Item(31, 5).Value = "  -0.06%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.05"
$ws.Cells.Item(32, 5).Value = "  +3.99%  "

$ws.Cells.Item(33, 4).Value = "1.76"
$ws.Cells.Item(33, 5).Value = "  -3.37%  "

$ws.Cells.Item(34, 4).Value = "0.0₃0820"
$ws.Cells.Item(34, 5).Value = "  +1.41%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.22"
$ws.Cells.Item(35, 5).Value = "  +6.31%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "168.74"
$ws.Cells.Item(36, 5).Value = "  -3.73%  "

$ws.Cells.Item(37, 5).Value = "  +1.33%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.00"
$ws.Cells.Item(38, 5).Value = "  -0.05%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.94"
$ws.Cells.Item(39, 5).Value = "  +6.89%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "19.12"
$ws.Cells.Item(40, 5).Value = "  +0.19%  "

$ws.Cells.Item(41, 5).Value = "  +0.03%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "169.15"
$ws.Cells.Item(42, 5).Value = "  -0.55%  "

$ws.Cells.Item(43, 5).Value = "  +1.07%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "22.22"
$ws.Cells.Item(44, 5).Value = "  -0.94%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0572"
$ws.Cells.Item(45, 5).Value = "  +2.84%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.633"
$ws.Cells.Item(46, 5).Value = "  +0.26%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0247"
$ws.Cells.Item(47, 5).Value = "  +3.54%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0962"
$ws.Cells.Item(48, 5).Value = "  +0.13%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "18.85"
$ws.Cells.Item(49, 5).Value = "  +0.38%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.89"
$ws.Cells.Item(50, 5).Value = "  +10.49%  "

$ws.Cells.Item(51, 2).Value = "TheGraph"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.177"
$ws.Cells.Item(51, 5).Value = "  +2.04%  "
